# Global_Links_to_CoinData.xlsx - add row 15 (sheet row 16) to the
# "Links_to_Data" table: a new numismatic marketplace link.
#
# & - №71036340 от 22.09.2024 https://meshok.net/

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Links_to_Data")
$ws.Select()

$newUrl = "https://www.coins-market.ru/"

# № column - next sequential index
$ws.Range("A16").Value = 15

# Links: column - URL text + live hyperlink
$ws.Hyperlinks.Add($ws.Range("B16"), $newUrl)
# Hyperlinks.Add() re-applies formatting (font/fill/border) on the anchor
# cell; restore the row's existing look (same as B14/B15) by copying the
# format from the cell directly above instead of leaving the default look.
$ws.Range("B15").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B16").Value = $newUrl

# Description: column - reuse the existing "Евро с описанием и тиражом"
# text/style used by the other Euro-catalogue rows (matches C14/C15).
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C16").Value = "Евро с описанием и тиражом"

# Leave the active selection on the newly filled row, as in the saved file.
$ws.Range("A16").Select()

# Best-effort: restore the window position recorded in the saved file.
$aw = $excel.ActiveWindow
if ($aw) {
    $aw.Left = 4250
    $aw.Top = 780
}
